$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.736.74'
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").Value = '''1.878.35'
$ws.Range("E3").Value = '  +1.42%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''331.51'
$ws.Range("E5").Value = '  +3.36%  '
$ws.Range("D6").Value = '''1.003'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '''0.4708'
$ws.Range("E7").Value = '  +5.02%  '
$ws.Range("D8").Value = '''0.3960'
$ws.Range("E8").Value = '  +3.05%  '
$ws.Range("D9").Value = '''47.91'
$ws.Range("E9").Value = '  -0.78%  '
$ws.Range("D10").Value = '''0.08029'
$ws.Range("E10").Value = '  +2.63%  '
$ws.Range("D11").Value = '''1.023'
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").Value = '''21.86'
$ws.Range("D13").Value = '''1.869.50'
$ws.Range("E13").Value = '  +2.94%  '
$ws.Range("D14").Value = '''5.968'
$ws.Range("E14").Value = '  +1.94%  '
$ws.Range("D15").Value = '''7.155'
$ws.Range("E15").Value = '  +0.82%  '
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '''87.15'
$ws.Range("E17").Value = '  +1.47%  '
$ws.Range("D18").Value = '''0.00001045'
$ws.Range("E18").Value = '  +2.12%  '
$ws.Range("D19").Value = '''0.06611'
$ws.Range("E19").Value = '  +1.77%  '
$ws.Range("D20").Value = '''17.29'
$ws.Range("E20").Value = '  +1.40%  '
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").Value = '''27.762.83'
$ws.Range("E22").Value = '  +1.31%  '
$ws.Range("D23").Value = '''5.510'
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("D24").Value = '''11.03'
$ws.Range("E24").Value = '  +2.27%  '
$ws.Range("D25").Value = '''2.297'
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("D26").Value = '''2.109.83'
$ws.Range("E26").Value = '  +3.28%  '
$ws.Range("D27").Value = '''156.52'
$ws.Range("E27").Value = '  +3.92%  '
$ws.Range("D28").Value = '''20.22'
$ws.Range("E28").Value = '  +4.76%  '
$ws.Range("D29").Value = '''2.092'
$ws.Range("E29").Value = '  +3.23%  '
$ws.Range("E30").Value = '  +2.20%  '
$ws.Range("D31").Value = '''122.56'
$ws.Range("E31").Value = '  +2.29%  '
$ws.Range("D32").Value = '''0.9703'
$ws.Range("E32").Value = '  +4.92%  '
$ws.Range("D33").Value = '''0.09552'
$ws.Range("E33").Value = '  +2.10%  '
$ws.Range("D34").Value = '''1.453'
$ws.Range("E34").Value = '  -3.05%  '
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("D36").Value = '''5.302'
$ws.Range("E36").Value = '  +1.43%  '
$ws.Range("D37").Value = '''0.06115'
$ws.Range("E37").Value = '  +2.99%  '
$ws.Range("D38").Value = '''0.02263'
$ws.Range("E38").Value = '  +2.05%  '
$ws.Range("D39").Value = '''1.232'
$ws.Range("E39").Value = '  +0.67%  '
$ws.Range("D40").Value = '''8.166'
$ws.Range("E40").Value = '  -1.30%  '
$ws.Range("D41").Value = '''0.5999'
$ws.Range("E41").Value = '  +1.85%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = '''0.1907'
$ws.Range("E43").Value = '  +3.58%  '
$ws.Range("E44").Value = '  -0.41%  '
$ws.Range("D45").Value = '''1.251'
$ws.Range("E45").Value = '  -2.19%  '
$ws.Range("D46").Value = '''0.5684'
$ws.Range("E46").Value = '  +1.35%  '
$ws.Range("D47").Value = '''12.22'
$ws.Range("E47").Value = '  +0.22%  '
$ws.Range("D48").Value = '''3.402'
$ws.Range("E48").Value = '  +1.52%  '
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("E51").Value = '  +9.58%  '
